# Applies the 2025-01-26 Sunday -> 2025-01-27 Monday dated worksheet update:
# the header date and all 25 practice problems/answers in the table are
# replaced with the new day's values.

$d = $word.ActiveDocument

$replacements = @(
    @("2025-01-26 Sunday", "2025-01-27 Monday"),
    @("195×6=1170", "179×2=358"),
    @("139×9=1251", "275×7=1925"),
    @("493×4=1972", "483×2=966"),
    @("293×4=1172", "259×4=1036"),
    @("407×9=3663", "867×8=6936"),
    @("458×4=1832", "952×2=1904"),
    @("450×6=2700", "766×6=4596"),
    @("355×2=710", "878×7=6146"),
    @("869×8=6952", "617×8=4936"),
    @("748×8=5984", "415×7=2905"),
    @("621×6=3726", "494×4=1976"),
    @("298×6=1788", "481×9=4329"),
    @("416×7=2912", "998×6=5988"),
    @("725×7=5075", "939×7=6573"),
    @("684×8=5472", "608×6=3648"),
    @("172×8=1376", "545×3=1635"),
    @("135×3=405", "935×7=6545"),
    @("740×8=5920", "345×2=690"),
    @("137×9=1233", "382×4=1528"),
    @("670×3=2010", "185×9=1665"),
    @("981×3=2943", "620×6=3720"),
    @("143×3=429", "160×5=800"),
    @("619×6=3714", "257×9=2313"),
    @("830×7=5810", "919×2=1838"),
    @("442×5=2210", "188×2=376")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
